$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as text so numeric-looking values
# (e.g. "1.001", "0.000007621") are not auto-converted to numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '30.442.80'
$ws.Range("E2").Value = '  -0.85%  '
$ws.Range("D3").Value = '1.917.68'
$ws.Range("E3").Value = '  +2.15%  '
$ws.Range("D4").Value = '1.000'
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").Value = '241.31'
$ws.Range("E5").Value = '  +1.81%  '
$ws.Range("D6").Value = '1.001'
$ws.Range("E6").Value = '  -0.01%  '
$ws.Range("D7").Value = '0.4707'
$ws.Range("E7").Value = '  -0.55%  '
$ws.Range("D8").Value = '0.2851'
$ws.Range("E8").Value = '  +1.36%  '
$ws.Range("D9").Value = '0.06797'
$ws.Range("E9").Value = '  +4.72%  '
$ws.Range("D10").Value = '106.77'
$ws.Range("E10").Value = '  +12.24%  '
$ws.Range("D11").Value = '18.27'
$ws.Range("E11").Value = '  -1.11%  '
$ws.Range("D12").Value = '1.895.89'
$ws.Range("E12").Value = '  +0.91%  '
$ws.Range("D13").Value = '0.07685'
$ws.Range("E13").Value = '  +1.67%  '
$ws.Range("D14").Value = '5.207'
$ws.Range("E14").Value = '  +2.92%  '
$ws.Range("D15").Value = '0.6567'
$ws.Range("E15").Value = '  +1.58%  '
$ws.Range("D16").Value = '289.08'
$ws.Range("E16").Value = '  -3.37%  '
$ws.Range("D17").Value = '30.441.29'
$ws.Range("E17").Value = '  -0.80%  '
$ws.Range("D18").Value = '0.000007621'
$ws.Range("E18").Value = '  +1.54%  '
$ws.Range("D19").Value = '1.001'
$ws.Range("E19").Value = '  +0.08%  '
$ws.Range("D20").Value = '12.95'
$ws.Range("E20").Value = '  -0.48%  '
$ws.Range("D21").Value = '2.156.57'
$ws.Range("E21").Value = '  +1.54%  '
$ws.Range("D22").Value = '1.000'
$ws.Range("E22").Value = '  -0.06%  '
$ws.Range("D23").Value = '5.223'
$ws.Range("E23").Value = '  +2.01%  '
$ws.Range("D24").Value = '6.196'
$ws.Range("E24").Value = '  +1.19%  '
$ws.Range("D25").Value = '168.22'
$ws.Range("E25").Value = '  -0.25%  '
$ws.Range("D26").Value = '9.299'
$ws.Range("E26").Value = '  +1.30%  '
$ws.Range("D27").Value = '21.51'
$ws.Range("E27").Value = '  +10.13%  '
$ws.Range("D28").Value = '2.073'
$ws.Range("E28").Value = '  +7.10%  '
$ws.Range("D29").Value = '0.1068'
$ws.Range("E29").Value = '  +1.22%  '
$ws.Range("D30").Value = '1.371'
$ws.Range("E30").Value = '  +1.26%  '
$ws.Range("D31").Value = '4.159'
$ws.Range("E31").Value = '  +0.51%  '
$ws.Range("D32").Value = '3.966'
$ws.Range("E32").Value = '  +0.91%  '
$ws.Range("D33").Value = '0.05051'
$ws.Range("E33").Value = '  +0.49%  '
$ws.Range("D34").Value = '0.7431'
$ws.Range("E34").Value = '  +3.74%  '
$ws.Range("D35").Value = '1.153'
$ws.Range("E35").Value = '  -0.90%  '
$ws.Range("D36").Value = '0.02093'
$ws.Range("E36").Value = '  +9.74%  '
$ws.Range("E37").Value = '  +1.30%  '
$ws.Range("D38").Value = '2.680'
$ws.Range("E38").Value = '  -0.88%  '
$ws.Range("D39").Value = '2.051'
$ws.Range("E39").Value = '  +0.56%  '
$ws.Range("D40").Value = '108.89'
$ws.Range("E40").Value = '  +1.86%  '
$ws.Range("D41").Value = '0.8690'
$ws.Range("E41").Value = '  -2.64%  '
$ws.Range("D42").Value = '5.859'
$ws.Range("E42").Value = '  +5.22%  '
$ws.Range("B43").Value = 'TheSandbox'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D43").Value = '0.4256'
$ws.Range("E43").Value = '  +2.18%  '
$ws.Range("B44").Value = 'PaxDollar'
$ws.Range("C44").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D44").Value = '1.001'
$ws.Range("E44").Value = '  +0.03%  '
$ws.Range("D45").Value = '67.55'
$ws.Range("E45").Value = '  +4.59%  '
$ws.Range("D46").Value = '50.42'
$ws.Range("E46").Value = '  +18.76%  '
$ws.Range("D47").Value = '7.173'
$ws.Range("E47").Value = '  -1.49%  '
$ws.Range("D48").Value = '9.256'
$ws.Range("E48").Value = '  +3.58%  '
$ws.Range("D49").Value = '0.1210'
$ws.Range("E49").Value = '  -0.17%  '
$ws.Range("D50").Value = '34.85'
$ws.Range("E50").Value = '  +1.16%  '
$ws.Range("D51").Value = '0.3914'
$ws.Range("E51").Value = '  +3.30%  '

# Restore the default (unstyled) look for column D so no stray cell-level
# style is left behind, matching the original workbook formatting.
$ws.Range("D2:D51").Style = "Normal"
